$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume cells are stored as text, matching the source data
# (these are textual display strings, not numeric values, e.g. "42.417.58")
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '42.417.58'
$ws.Range("E2").Value = '  +1.72%  '
$ws.Range("D3").Value = '2.288.58'
$ws.Range("E3").Value = '  +0.97%  '
$ws.Range("D5").Value = '156.98'
$ws.Range("E5").Value = '  +15,594.98%  '
$ws.Range("D6").Value = '307.69'
$ws.Range("E6").Value = '  +1.15%  '
$ws.Range("D7").Value = '95.65'
$ws.Range("E8").Value = '  +0.37%  '
$ws.Range("E9").Value = '  -0.06%  '
$ws.Range("D10").Value = '0.495'
$ws.Range("E10").Value = '  +2.98%  '
$ws.Range("D11").Value = '35.87'
$ws.Range("E11").Value = '  +11.17%  '
$ws.Range("D12").Value = '0.0803'
$ws.Range("E12").Value = '  +1.04%  '
$ws.Range("E13").Value = '  -2.05%  '
$ws.Range("D14").Value = '6.73'
$ws.Range("E14").Value = '  +2.32%  '
$ws.Range("D15").Value = '2.641.49'
$ws.Range("E15").Value = '  +0.93%  '
$ws.Range("D16").Value = '14.49'
$ws.Range("E16").Value = '  +2.21%  '
$ws.Range("D17").Value = '2.290.33'
$ws.Range("E17").Value = '  +1.20%  '
$ws.Range("D18").Value = '0.799'
$ws.Range("E18").Value = '  +5.11%  '
$ws.Range("D19").Value = '42.336.97'
$ws.Range("E19").Value = '  +1.67%  '
$ws.Range("D20").Value = '12.71'
$ws.Range("E20").Value = '  +3.32%  '
$ws.Range("D21").Value = '0.0₃0917'
$ws.Range("E21").Value = '  +1.68%  '
$ws.Range("D22").Value = '6.01'
$ws.Range("E22").Value = '  +1.88%  '
$ws.Range("D23").Value = '68.02'
$ws.Range("E23").Value = '  +2.09%  '
$ws.Range("D24").Value = '242.84'
$ws.Range("E24").Value = '  +0.76%  '
$ws.Range("D25").Value = '2.62'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("E26").Value = '  +2.12%  '
$ws.Range("E27").Value = '  -0.11%  '
$ws.Range("D28").Value = '24.07'
$ws.Range("E28").Value = '  -0.26%  '
$ws.Range("D29").Value = '35.96'
$ws.Range("E29").Value = '  +4.27%  '
$ws.Range("D30").Value = '9.58'
$ws.Range("E30").Value = '  +0.89%  '
$ws.Range("E31").Value = '  -2.44%  '
$ws.Range("D32").Value = '161.53'
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("D33").Value = '5.34'
$ws.Range("E33").Value = '  +3.82%  '
$ws.Range("D34").Value = '0.999'
$ws.Range("E34").Value = '  +0.03%  '
$ws.Range("D35").Value = '0.0755'
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("D36").Value = '3.10'
$ws.Range("E36").Value = '  +2.66%  '
$ws.Range("E37").Value = '  +4.49%  '
$ws.Range("D38").Value = '17.25'
$ws.Range("E38").Value = '  +4.09%  '
$ws.Range("E39").Value = '  -0.15%  '
$ws.Range("E40").Value = '  +2.49%  '
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").Value = '4.16'
$ws.Range("E42").Value = '  +6.92%  '
$ws.Range("D43").Value = '2.010.29'
$ws.Range("E43").Value = '  -2.44%  '
$ws.Range("D44").Value = '19.63'
$ws.Range("E44").Value = '  +0.87%  '
$ws.Range("E45").Value = '  +11.06%  '
$ws.Range("D46").Value = '0.0285'
$ws.Range("E46").Value = '  +2.69%  '
$ws.Range("D47").Value = '10.15'
$ws.Range("E47").Value = '  +0.26%  '
$ws.Range("E48").Value = '  +4.57%  '
$ws.Range("E49").Value = '  +2.61%  '
$ws.Range("D50").Value = '53.33'
$ws.Range("E50").Value = '  +3.40%  '
$ws.Range("D51").Value = '72.90'
$ws.Range("E51").Value = '  -0.17%  '
